$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: the manager's first reply got new text + new timestamp
$ws.Range("E5").Value = "Yes I can see this!"
$ws.Range("H5").Value = "2025-04-23T15:26:51.459171"

# Remove the duplicated row 7 (it duplicated row 4/row 5's content a second time)
$ws.Rows.Item(7).Delete()

# Row 6 (now the last row) becomes a brand-new manager reply "How about you?"
# Copy the "2" text value (Message ID) from an existing text cell so it keeps
# its Text data type (not a literal number) without altering cell styling.
$ws.Range("A3").Copy($ws.Range("D6"))
$ws.Range("E6").Value = "How about you?"
$ws.Range("F6").Value = "S5678901G"
$ws.Range("G6").Value = "Manager"
$ws.Range("H6").Value = "2025-04-23T15:27:06.254456"
